$p = $ppt.ActivePresentation
$newSlide = $p.Slides.Add(5, 2)
$left = 838200 / 12700.0
$top = 3269774 / 12700.0
$width = 10515600 / 12700.0
$height = 1463040 / 12700.0
$tbl = $newSlide.Shapes.AddTable(4, 4, $left, $top, $width, $height)
$newSlide.Shapes.Item(1).Delete()
$newSlide.Shapes.Item(1).Delete()
$tbl.Name = "Content Placeholder 3"

$cell = $tbl.Table.Cell(1,1)
$tr = $cell.Shape.TextFrame.TextRange
$tr.Text = "Activity"
Write-Host "LanguageID before: $($tr.LanguageID)"
$tr.LanguageID = 2057
Write-Host "LanguageID after: $($tr.LanguageID)"
$tr.Font.Bold = -1
